# "Generate Report for Handback"
# Renames the locale "ru-ru" to "fr-fr" throughout the workbook (sheet name,
# Overview header cell, both table names/columns) and refreshes the
# Correspond Handoff / Handback datetime stamps.

$wb = $excel.ActiveWorkbook

$wsLocale = $wb.Worksheets.Item("ru-ru")
$wsOverview = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------------------
# 1. Update the Overview sheet's header cell (B1) that names the locale.
# ---------------------------------------------------------------------------
$wsOverview.Range("B1").Value = "fr-fr"

# ---------------------------------------------------------------------------
# 2. Fix up the Overview table (table id=2) so its second column is named
#    "fr-fr" too. The ListColumns.Item(n).Name setter is a no-op in this
#    runtime for header-less tables (headerRowCount=0), so the table is
#    unlisted and re-created against the (already updated) B1 cell, then the
#    lost attributes (name/displayName/style/header setting) are restored.
# ---------------------------------------------------------------------------
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Unlist()

$newLoOverview = $wsOverview.ListObjects.Add(1, $wsOverview.Range("A1:C1"), $null, $false)
$newLoOverview.ShowHeaders = $false
# Re-creating + toggling ShowHeaders shifts the table ref down one row;
# resize it back to the original single-row range.
$newLoOverview.Resize($wsOverview.Range("A1:C1"))
$newLoOverview.Name = "Overview"
$newLoOverview.TableStyle = "TableStyleMedium9"
# Re-creating the table stamped a placeholder "Column3" literal into C1
# (which was blank before); clear it back out.
$wsOverview.Range("C1").ClearContents()

# ---------------------------------------------------------------------------
# 3. Rename the locale worksheet itself.
# ---------------------------------------------------------------------------
$wsLocale.Name = "fr-fr"

# ---------------------------------------------------------------------------
# 4. Rename the locale sheet's table (table id=1) to match.
# ---------------------------------------------------------------------------
$loLocale = $wsLocale.ListObjects.Item(1)
$loLocale.Name = "fr-fr"

# ---------------------------------------------------------------------------
# 5. Refresh the handoff / handback timestamps shown on the locale sheet.
#    Every "Correspond Handoff Datetime" cell (column E) now reads the same
#    refreshed timestamp, and every "Correspond Handback DateTime" cell
#    (column H) reads the new, shared, handback timestamp.
# ---------------------------------------------------------------------------
$wsLocale.Range("E2").Value = "2016-03-11 01:02:43"
$wsLocale.Range("E3").Value = "2016-03-11 01:02:43"
$wsLocale.Range("E4").Value = "2016-03-11 01:02:43"
$wsLocale.Range("E5").Value = "2016-03-11 01:02:43"

$wsLocale.Range("H2").Value = "2016-03-17 17:12:56"
$wsLocale.Range("H3").Value = "2016-03-17 17:12:56"
$wsLocale.Range("H4").Value = "2016-03-17 17:12:56"
$wsLocale.Range("H5").Value = "2016-03-17 17:12:56"
